# Update "想去人数" (interest count) figures in column F across all sheets,
# matching the gh-pages data refresh at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1746
$ws.Range("F4").Value = 159
$ws.Range("F5").Value = 442
$ws.Range("F6").Value = 816
$ws.Range("F7").Value = 245
$ws.Range("F8").Value = 1177
$ws.Range("F9").Value = 334
$ws.Range("F11").Value = 876
$ws.Range("F12").Value = 678
$ws.Range("F14").Value = 505
$ws.Range("F17").Value = 170
$ws.Range("F18").Value = 2911
$ws.Range("F19").Value = 2616
$ws.Range("F21").Value = 28
$ws.Range("F23").Value = 315
$ws.Range("F26").Value = 5258
$ws.Range("F28").Value = 977
$ws.Range("F31").Value = 300
$ws.Range("F32").Value = 1085

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 1115
$ws.Range("F14").Value = 608
$ws.Range("F19").Value = 42
$ws.Range("F24").Value = 312
$ws.Range("F25").Value = 277
$ws.Range("F26").Value = 3915
$ws.Range("F33").Value = 160

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 2446
$ws.Range("F6").Value = 1031
$ws.Range("F9").Value = 1313
$ws.Range("F11").Value = 99

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 2446
$ws.Range("F5").Value = 1746
$ws.Range("F6").Value = 1031
$ws.Range("F7").Value = 1313
$ws.Range("F9").Value = 99
$ws.Range("F10").Value = 159
$ws.Range("F11").Value = 442
$ws.Range("F12").Value = 816
$ws.Range("F13").Value = 245
$ws.Range("F14").Value = 1177
$ws.Range("F15").Value = 334
$ws.Range("F16").Value = 876
$ws.Range("F17").Value = 678
$ws.Range("F18").Value = 1115
$ws.Range("F19").Value = 1115
$ws.Range("F20").Value = 505
$ws.Range("F22").Value = 170
$ws.Range("F23").Value = 2911
$ws.Range("F24").Value = 2616
$ws.Range("F26").Value = 28
$ws.Range("F27").Value = 315
$ws.Range("F30").Value = 5258
$ws.Range("F32").Value = 977
$ws.Range("F33").Value = 608
$ws.Range("F34").Value = 608
$ws.Range("F38").Value = 300
$ws.Range("F41").Value = 42
$ws.Range("F44").Value = 312
$ws.Range("F45").Value = 312
$ws.Range("F46").Value = 277
$ws.Range("F47").Value = 1085
$ws.Range("F49").Value = 160

$wb.Save()
